$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 21: new blade entry (opposite direction / reverse blades) - entered first
$ws.Range("G21").Value = "http://www.hobbyexpress.com/gemfan_11x4.7_reverse_carbon_filled_1041840_prd1.htm?pSearchQueryId=4684805"
$ws.Range("B21").Value = "Blades opposite direction"
$ws.Range("D21").Value = 5
$ws.Range("E21").Value = 3

# Row 20: update blade cost + source link (normal rotation blades, now sourced from Hobby Express)
$ws.Range("A20").ClearContents()
$ws.Range("E20").Value = 3
$ws.Range("G20").Value = "http://www.hobbyexpress.com/gemfan_11x4.7_normal_carbon_filled_1041839_prd1.htm?pSearchQueryId=4684805"

$ws.Range("B15").Select()
